# Applies the scheduled-runner profit/price refresh to the Leve profit sheets.
# For each (sheet, row) the currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are updated to the freshly computed values; where a column
# has no longer any value (e.g. a profit no longer computable), its cell is
# cleared instead of being set to a new number.
$wb = $excel.ActiveWorkbook

# ALC!40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11508.565
$ws.Range("I40").Value = 5099.6665
$ws.Range("J40").Value = 12469.9
$ws.Range("K40").Value = 5099.6665
$ws.Range("L40").Value = 12469.9
$ws.Range("M40").Value = -4924.6665
$ws.Range("N40").Value = -12819.9

# ALC!117
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 81475.25
$ws.Range("J117").Value = 81475.25
$ws.Range("L117").Value = 81475.25
$ws.Range("N117").Value = -90653.25

# ALC!120
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 49493.75
$ws.Range("J120").Value = 49493.75
$ws.Range("L120").Value = 49493.75
$ws.Range("N120").Value = -59169.75

# ALC!132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1828.6818
$ws.Range("I132").Value = 1612.8975
$ws.Range("K132").Value = 4838.6925
$ws.Range("M132").Value = -2308.6925

# ALC!138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2039.2084
$ws.Range("I138").Value = 1243.3914
$ws.Range("J138").Value = 2771.36
$ws.Range("K138").Value = 3730.1742
$ws.Range("L138").Value = 8314.08
$ws.Range("M138").Value = 1409.8258
$ws.Range("N138").Value = -18594.08

# ARM!23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 129999
$ws.Range("J23").Value = 129999
$ws.Range("L23").Value = 129999
$ws.Range("N23").Value = -130517

# ARM!44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 96144
$ws.Range("I44").Value = 70000
$ws.Range("K44").Value = 70000
$ws.Range("M44").Value = -69512

# ARM!45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15627056
$ws.Range("I45").Value = 2438.8333
$ws.Range("K45").Value = 2438.8333
$ws.Range("M45").Value = -2061.8333

# ARM!104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 34733.668
$ws.Range("J104").Value = 34733.668
$ws.Range("L104").Value = 34733.668
$ws.Range("N104").Value = -41721.668

# ARM!117
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 70706.57000000001
$ws.Range("J117").Value = 70706.57000000001
$ws.Range("L117").Value = 70706.57000000001
$ws.Range("N117").Value = -79884.57000000001

# ARM!132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2379.7856
$ws.Range("J132").Value = 3599.6
$ws.Range("L132").Value = 10798.8
$ws.Range("N132").Value = -15858.8

# BSM!6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 5515.3335
$ws.Range("J6").Value = 6500
$ws.Range("L6").Value = 6500
$ws.Range("N6").Value = -6726

# BSM!35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 64443.668
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 64443.668
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 64443.668
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -65063.668

# BSM!110
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 75484
$ws.Range("J110").Value = 75484
$ws.Range("L110").Value = 75484
$ws.Range("N110").Value = -83664

# BSM!122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 93644.42999999999
$ws.Range("J122").Value = 93644.42999999999
$ws.Range("L122").Value = 93644.42999999999
$ws.Range("N122").Value = -103444.43

# BSM!134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3381.2942
$ws.Range("J134").Value = 5246.8887
$ws.Range("L134").Value = 15740.6661
$ws.Range("N134").Value = -20810.6661

# CRP!9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 33314.715
$ws.Range("J9").Value = 33314.715
$ws.Range("L9").Value = 33314.715
$ws.Range("N9").Value = -33650.715

# CRP!31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3139.037
$ws.Range("I31").Value = 2333.1667
$ws.Range("K31").Value = 2333.1667
$ws.Range("M31").Value = -2038.1667

# CRP!34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3139.037
$ws.Range("I34").Value = 2333.1667
$ws.Range("K34").Value = 2333.1667
$ws.Range("M34").Value = -2131.1667

# CRP!108
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 70801.39999999999
$ws.Range("J108").Value = 70801.39999999999
$ws.Range("L108").Value = 70801.39999999999
$ws.Range("N108").Value = -78481.39999999999

# CRP!117
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 44081.168
$ws.Range("J117").Value = 44081.168
$ws.Range("L117").Value = 44081.168
$ws.Range("N117").Value = -53259.168

# CRP!119
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 63949.43
$ws.Range("J119").Value = 63949.43
$ws.Range("L119").Value = 63949.43
$ws.Range("N119").Value = -73625.42999999999

# CRP!134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 25821.658
$ws.Range("I134").Value = 2954.6667
$ws.Range("K134").Value = 8864.000100000001
$ws.Range("M134").Value = -6329.000100000001

# CRP!138
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 109369.71
$ws.Range("J138").Value = 116646.5
$ws.Range("L138").Value = 116646.5
$ws.Range("N138").Value = -126926.5

# CUL!19
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 20
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 60
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 114
$ws.Range("N19").ClearContents()

# CUL!34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 232
$ws.Range("J34").Value = 75
$ws.Range("L34").Value = 225
$ws.Range("N34").Value = -393

# CUL!52
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2339
$ws.Range("J52").Value = 3397.5
$ws.Range("L52").Value = 10192.5
$ws.Range("N52").Value = -10724.5

# GSM!102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3156.889
$ws.Range("I102").Value = 3214
$ws.Range("K102").Value = 3214
$ws.Range("M102").Value = -1592

# GSM!108
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 57228.168
$ws.Range("J108").Value = 57228.168
$ws.Range("L108").Value = 57228.168
$ws.Range("N108").Value = -64908.168

# GSM!132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3952.1714
$ws.Range("I132").Value = 3461.08
$ws.Range("K132").Value = 10383.24
$ws.Range("M132").Value = -7853.24

# LTW!55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4772.5835
$ws.Range("I55").Value = 1142.2609
$ws.Range("K55").Value = 1142.2609
$ws.Range("M55").Value = -969.2609

# LTW!117
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 30565.111
$ws.Range("J117").Value = 30565.111
$ws.Range("L117").Value = 30565.111
$ws.Range("N117").Value = -39743.111

# LTW!118
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 60143.11
$ws.Range("J118").Value = 60143.11
$ws.Range("L118").Value = 60143.11
$ws.Range("N118").Value = -63457.11

# WVR!54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33783.332
$ws.Range("I54").Value = 28675
$ws.Range("K54").Value = 28675
$ws.Range("M54").Value = -28155

# WVR!113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1764.6666
$ws.Range("I113").Value = 1911.7142
$ws.Range("K113").Value = 5735.142599999999
$ws.Range("M113").Value = -3565.142599999999

# WVR!121
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 53885.332
$ws.Range("J121").Value = 53885.332
$ws.Range("L121").Value = 53885.332
$ws.Range("N121").Value = -57379.332

# WVR!132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2559793.5
$ws.Range("I132").Value = 2348.4285
$ws.Range("K132").Value = 7045.2855
$ws.Range("M132").Value = -4515.2855
